# Generate Report for Handoff
# Updates the localization-status workbook to reflect that translations are
# now "Ready for handoff" (previously "In Translation"), with refreshed
# handoff timestamps for the Overview sheet and each language sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn (B2) and de-de (C2), plus the
# combined "Latest Handoff Date" column (D2).
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-32-13 00:32:32"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2).
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-13 00:32:29"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2).
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-13 00:32:32"
